$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Joby Aviation / JOBY)
$ws.Range("D2").Value = 15.48
$ws.Range("E2").Value = 55.6
$ws.Range("F2").Value = 9.67
$ws.Range("N2").Value = 52.47848103381103

# Row 3 (Archer Aviation / ACHR)
$ws.Range("D3").Value = 8.640000000000001
$ws.Range("E3").Value = 58.5
$ws.Range("F3").Value = 15.38
$ws.Range("H3").Value = 40
$ws.Range("I3").Value = 53
$ws.Range("K3").Value = 54.9
$ws.Range("N3").Value = 52.47848103381103
